# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 2.9
$ws.Range("I2").Value = 2.75
$ws.Range("L2").Value = 3.6
$ws.Range("X2").Value = 12
$ws.Range("AF2").Value = 67
$ws.Range("AI2").Value = 12
$ws.Range("AL2").Value = 29
$ws.Range("H3").Value = 5.6
$ws.Range("I3").Value = 11.25
$ws.Range("K3").Value = 2.65
$ws.Range("L3").Value = 8.5
$ws.Range("S3").Value = 1.28
$ws.Range("T3").Value = 3.5
$ws.Range("U3").Value = 2.05
$ws.Range("V3").Value = 1.7
$ws.Range("X3").Value = 6.6
$ws.Range("Y3").Value = 9.75
$ws.Range("Z3").Value = 7.8
$ws.Range("AA3").Value = 11.5
$ws.Range("AB3").Value = 32
$ws.Range("AI3").Value = 100
$ws.Range("AJ3").Value = 37
$ws.Range("AK3").Value = 400
$ws.Range("AL3").Value = 150
$ws.Range("AN3").Value = 3.15
$ws.Range("AT3").Value = 3.5
$ws.Range("AW3").Value = 11
$ws.Range("AX3").Value = 60
$ws.Range("AY3").Value = 50
$ws.Range("BA3").Value = 400
$ws.Range("M4").Value = 1.03
$ws.Range("O4").Value = 1.17
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 4.2
$ws.Range("I5").Value = 1.48
$ws.Range("K5").Value = 2.4
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("O5").Value = 1.18
$ws.Range("P5").Value = 4.5
$ws.Range("Q5").Value = 1.65
$ws.Range("R5").Value = 2.2
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 1.91
$ws.Range("W5").Value = 21
$ws.Range("AB5").Value = 51
$ws.Range("AD5").Value = 8
$ws.Range("AG5").Value = 201
$ws.Range("AZ5").Value = 21
$ws.Range("BB5").Value = 101
$ws.Range("U6").Value = 1.75
$ws.Range("U7").Value = 1.67
$ws.Range("Q8").Value = 2.08
$ws.Range("R8").Value = 1.73
$ws.Range("U8").Value = 1.83
$ws.Range("V8").Value = 1.83
$ws.Range("S9").Value = 1.29
$ws.Range("T9").Value = 3.5
$ws.Range("X9").Value = 8.5
$ws.Range("AS9").Value = 101
$ws.Range("AT9").Value = 3.5
$ws.Range("BA9").Value = 101
